# Natmi following Dr Hou advice
#
# The ligand-receptor table is recomputed: "Sending cluster" now also
# includes the "sCs" cluster (in addition to "FAPs"), each paired against
# all four target clusters (ECs, FAPs, M2, sCs). This grows the data block
# from 4 rows (A2:T5) to 8 rows (A2:T9), and refreshes every numeric
# statistic with the newly recalculated values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: Sending cluster, Ligand symbol, Receptor symbol, Target cluster,
# then the 16 numeric statistic columns (E..T).
$data = @(
    @("FAPs","Lgi3","Flot1","ECs", 3,1,1.535771,4.607313,0.8955152254457543,0.8955152254457543,3,1,23.18420866666666,69.55262599999999,0.3933108221010576,0.3933108221010577,35.60563532821534,320.450717953938,0.3522158295240836,0.3522158295240836),
    @("FAPs","Lgi3","Flot1","FAPs",3,1,1.535771,4.607313,0.8955152254457543,0.8955152254457543,3,1,13.52213766666667,40.566413,0.2293976541837692,0.2293976541837692,20.76690688647433,186.902161978269,0.2054290920031052,0.2054290920031053),
    @("FAPs","Lgi3","Flot1","M2",  3,1,1.535771,4.607313,0.8955152254457543,0.8955152254457543,3,1,14.10188666666667,42.30566,0.239232863963006,0.239232863963006,21.65726858795334,194.91541729158,0.2142366721058648,0.2142366721058648),
    @("FAPs","Lgi3","Flot1","sCs", 3,1,1.535771,4.607313,0.8955152254457543,0.8955152254457543,3,1,8.138043999999999,24.414132,0.1380586597521672,0.1380586597521672,12.498171971924,112.483547747316,0.1236336318127007,0.1236336318127007),
    @("sCs","Lgi3","Flot1","ECs",  1,0.3333333333333333,0.179187,0.537561,0.1044847745542456,0.1044847745542456,3,1,23.18420866666666,69.55262599999999,0.3933108221010576,0.3933108221010577,4.154308798353999,37.38877918518599,0.04109499257697401,0.04109499257697402),
    @("sCs","Lgi3","Flot1","FAPs", 1,0.3333333333333333,0.179187,0.537561,0.1044847745542456,0.1044847745542456,3,1,13.52213766666667,40.566413,0.2293976541837692,0.2293976541837692,2.422991282077,21.806921538693,0.02396856218066392,0.02396856218066393),
    @("sCs","Lgi3","Flot1","M2",   1,0.3333333333333333,0.179187,0.537561,0.1044847745542456,0.1044847745542456,3,1,14.10188666666667,42.30566,0.239232863963006,0.239232863963006,2.52687476614,22.74187289526,0.02499619185714119,0.0249961918571412),
    @("sCs","Lgi3","Flot1","sCs",  1,0.3333333333333333,0.179187,0.537561,0.1044847745542456,0.1044847745542456,3,1,8.138043999999999,24.414132,0.1380586597521672,0.1380586597521672,1.458231690228,13.124085212052,0.0144250279394665,0.0144250279394665)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $rowValues = $data[$i]
    for ($c = 1; $c -le $rowValues.Length; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowValues[$c - 1]
    }
}
